# Update "想去人数" (number of people interested) counts across the
# workbook's sheets. Each of the first three sheets (展览 / 演出 / 本地生活)
# holds its own category of events, and the fourth sheet (全部类型)
# aggregates the same rows again, so every value needs to be bumped in
# both places.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F30").Value = 4938
$wsExhibition.Range("F31").Value = 426
$wsExhibition.Range("F32").Value = 191

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F24").Value = 362
$wsShow.Range("F26").Value = 615
$wsShow.Range("F33").Value = 243

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F5").Value = 397

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 397
$wsAll.Range("F40").Value = 362
$wsAll.Range("F41").Value = 615
$wsAll.Range("F46").Value = 4938
$wsAll.Range("F48").Value = 426
$wsAll.Range("F49").Value = 191
